$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 113
$ws.Range("H113").Value = 10037
$ws.Range("I113").Value = 11500
$ws.Range("J113").Value = 9619
$ws.Range("K113").Value = 11500
$ws.Range("L113").Value = 9619
$ws.Range("M113").Value = -8246
$ws.Range("N113").Value = -16127
# Row 116
$ws.Range("H116").Value = 2687.6
$ws.Range("J116").Value = 3222
$ws.Range("L116").Value = 3222
$ws.Range("N116").Value = -10106
# Row 121
$ws.Range("H121").Value = 3936.375
$ws.Range("J121").Value = 3936.375
$ws.Range("L121").Value = 11809.125
$ws.Range("N121").Value = -15303.125
# Row 132
$ws.Range("H132").Value = 27795.21
$ws.Range("I132").Value = 1594.6471
$ws.Range("J132").Value = 250500
$ws.Range("K132").Value = 4783.9413
$ws.Range("L132").Value = 751500
$ws.Range("M132").Value = -2253.9413
$ws.Range("N132").Value = -756560
# Row 138
$ws.Range("H138").Value = 2610.5454
$ws.Range("I138").Value = 1491.1538
$ws.Range("J138").Value = 4227.4443
$ws.Range("K138").Value = 4473.4614
$ws.Range("L138").Value = 12682.3329
$ws.Range("M138").Value = 666.5385999999999
$ws.Range("N138").Value = -22962.3329

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2757.6667
$ws.Range("I2").Value = 2727.375
$ws.Range("K2").Value = 2727.375
$ws.Range("M2").Value = -2614.375
# Row 45
$ws.Range("H45").Value = 2282.4
$ws.Range("I45").Value = 1758.2222
$ws.Range("K45").Value = 1758.2222
$ws.Range("M45").Value = -1381.2222
# Row 74
$ws.Range("H74").Value = 2147.5
$ws.Range("I74").Value = 2197.353
$ws.Range("K74").Value = 2197.353
$ws.Range("M74").Value = -1323.353
# Row 77
$ws.Range("H77").Value = 2147.5
$ws.Range("I77").Value = 2197.353
$ws.Range("K77").Value = 10986.765
$ws.Range("M77").Value = -6618.764999999999
# Row 116
$ws.Range("H116").Value = 2757.6667
$ws.Range("I116").Value = 2727.375
$ws.Range("K116").Value = 2727.375
$ws.Range("M116").Value = -433.375
# Row 132
$ws.Range("H132").Value = 2123.0444
$ws.Range("I132").Value = 1382.5714
$ws.Range("J132").Value = 4714.7
$ws.Range("K132").Value = 4147.7142
$ws.Range("L132").Value = 14144.1
$ws.Range("M132").Value = -1617.7142
$ws.Range("N132").Value = -19204.1

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2757.6667
$ws.Range("I3").Value = 2727.375
$ws.Range("K3").Value = 2727.375
$ws.Range("M3").Value = -2613.375
# Row 64
$ws.Range("H64").Value = 1285.75
$ws.Range("I64").Value = 1199
$ws.Range("J64").Value = 1314.6666
$ws.Range("K64").Value = 1199
$ws.Range("L64").Value = 1314.6666
$ws.Range("M64").Value = -974
$ws.Range("N64").Value = -1764.6666
# Row 67
$ws.Range("H67").Value = 1285.75
$ws.Range("I67").Value = 1199
$ws.Range("J67").Value = 1314.6666
$ws.Range("K67").Value = 1199
$ws.Range("L67").Value = 1314.6666
$ws.Range("M67").Value = -419
$ws.Range("N67").Value = -2874.6666
# Row 86
$ws.Range("H86").Value = 2550.0588
$ws.Range("I86").Value = 1180.8462
$ws.Range("J86").Value = 7000
$ws.Range("K86").Value = 1180.8462
$ws.Range("L86").Value = 7000
$ws.Range("M86").Value = -57.84619999999995
$ws.Range("N86").Value = -9246
# Row 89
$ws.Range("H89").Value = 2550.0588
$ws.Range("I89").Value = 1180.8462
$ws.Range("J89").Value = 7000
$ws.Range("K89").Value = 5904.231
$ws.Range("L89").Value = 35000
$ws.Range("M89").Value = -288.2309999999998
$ws.Range("N89").Value = -46232
# Row 99
$ws.Range("H99").Value = 2850.125
$ws.Range("I99").Value = 2798.4
$ws.Range("J99").Value = 2936.3333
$ws.Range("K99").Value = 2798.4
$ws.Range("L99").Value = 2936.3333
$ws.Range("M99").Value = -1300.4
$ws.Range("N99").Value = -5932.3333
# Row 103
$ws.Range("H103").Value = 10885.667
$ws.Range("J103").Value = 10885.667
$ws.Range("L103").Value = 10885.667
$ws.Range("N103").Value = -13229.667
# Row 105
$ws.Range("H105").Value = 5099.9165
$ws.Range("I105").Value = 4591.6665
$ws.Range("K105").Value = 4591.6665
$ws.Range("M105").Value = -2844.6665

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 52
$ws.Range("H52").Value = 53220.8
$ws.Range("I52").Value = 60275
$ws.Range("J52").Value = 48518
$ws.Range("K52").Value = 60275
$ws.Range("L52").Value = 48518
$ws.Range("M52").Value = -59981
$ws.Range("N52").Value = -49106
# Row 99
$ws.Range("H99").Value = 8316762.5
$ws.Range("J99").Value = 15392054
$ws.Range("L99").Value = 15392054
$ws.Range("N99").Value = -15395050
# Row 126
$ws.Range("H126").Value = 8316762.5
$ws.Range("J126").Value = 15392054
$ws.Range("L126").Value = 46176162
$ws.Range("N126").Value = -46181102
# Row 132
$ws.Range("H132").Value = 2204.2307
$ws.Range("I132").Value = 1445.9546
$ws.Range("K132").Value = 4337.8638
$ws.Range("M132").Value = -1807.8638

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 3220.48
$ws.Range("I5").Value = 2119.6
$ws.Range("K5").Value = 6358.799999999999
$ws.Range("M5").Value = -6246.799999999999
# Row 57
$ws.Range("H57").Value = 202224.5
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 202224.5
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 606673.5
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -607791.5
# Row 131
$ws.Range("H131").Value = 3772.1
$ws.Range("I131").Value = 2617
$ws.Range("K131").Value = 7851
$ws.Range("M131").Value = -2811
# Row 135
$ws.Range("H135").Value = 3220.48
$ws.Range("I135").Value = 2119.6
$ws.Range("K135").Value = 19076.4
$ws.Range("M135").Value = -16541.4

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 25730150
$ws.Range("I80").Value = 226161
$ws.Range("J80").Value = 41670144
$ws.Range("K80").Value = 226161
$ws.Range("L80").Value = 41670144
$ws.Range("M80").Value = -225163
$ws.Range("N80").Value = -41672140
# Row 83
$ws.Range("H83").Value = 25730150
$ws.Range("I83").Value = 226161
$ws.Range("J83").Value = 41670144
$ws.Range("K83").Value = 1130805
$ws.Range("L83").Value = 208350720
$ws.Range("M83").Value = -1125813
$ws.Range("N83").Value = -208360704
# Row 102
$ws.Range("H102").Value = 1487.8096
$ws.Range("I102").Value = 1444.4736
$ws.Range("K102").Value = 1444.4736
$ws.Range("M102").Value = 177.5264
# Row 132
$ws.Range("H132").Value = 1812.3182
$ws.Range("I132").Value = 1554.2354
$ws.Range("J132").Value = 2689.8
$ws.Range("K132").Value = 4662.706200000001
$ws.Range("L132").Value = 8069.400000000001
$ws.Range("M132").Value = -2132.706200000001
$ws.Range("N132").Value = -13129.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 6660.1763
$ws.Range("I40").Value = 7356.5
$ws.Range("K40").Value = 7356.5
$ws.Range("M40").Value = -7220.5
# Row 46
$ws.Range("H46").Value = 783.375
$ws.Range("I46").Value = 744.6667
$ws.Range("J46").Value = 833.1429000000001
$ws.Range("K46").Value = 744.6667
$ws.Range("L46").Value = 833.1429000000001
$ws.Range("M46").Value = -556.6667
$ws.Range("N46").Value = -1209.1429
# Row 82
$ws.Range("H82").Value = 4264.4443
$ws.Range("I82").Value = 2298.3333
$ws.Range("J82").Value = 8196.666999999999
$ws.Range("K82").Value = 2298.3333
$ws.Range("L82").Value = 8196.666999999999
$ws.Range("M82").Value = -1937.3333
$ws.Range("N82").Value = -8918.666999999999
# Row 85
$ws.Range("H85").Value = 4264.4443
$ws.Range("I85").Value = 2298.3333
$ws.Range("J85").Value = 8196.666999999999
$ws.Range("K85").Value = 2298.3333
$ws.Range("L85").Value = 8196.666999999999
$ws.Range("M85").Value = -1050.3333
$ws.Range("N85").Value = -10692.667

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 10103892
$ws.Range("I81").Value = 3513.8572
$ws.Range("K81").Value = 7027.7144
$ws.Range("M81").Value = -5966.7144
# Row 84
$ws.Range("H84").Value = 10103892
$ws.Range("I84").Value = 3513.8572
$ws.Range("K84").Value = 35138.572
$ws.Range("M84").Value = -29834.572
# Row 130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
